$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 2040
$ws.Range("I29").Value = 120
$ws.Range("K29").Value = 360
$ws.Range("M29").Value = -79
$ws.Range("H38").Value = 631.44446
$ws.Range("I38").Value = 302.75
$ws.Range("J38").Value = 1288.8334
$ws.Range("K38").Value = 908.25
$ws.Range("L38").Value = 3866.5002
$ws.Range("M38").Value = -536.25
$ws.Range("N38").Value = -4610.5002
$ws.Range("H40").Value = 2817.4119
$ws.Range("I40").Value = 4557.8
$ws.Range("J40").Value = 2092.25
$ws.Range("K40").Value = 4557.8
$ws.Range("L40").Value = 2092.25
$ws.Range("M40").Value = -4382.8
$ws.Range("N40").Value = -2442.25
$ws.Range("H41").Value = 633
$ws.Range("I41").Value = 340.6154
$ws.Range("J41").Value = 1900
$ws.Range("K41").Value = 340.6154
$ws.Range("L41").Value = 1900
$ws.Range("M41").Value = 99.38459999999998
$ws.Range("N41").Value = -2780
$ws.Range("H58").Value = 68620
$ws.Range("I58").Value = 1325
$ws.Range("J58").Value = 93090.91
$ws.Range("K58").Value = 3975
$ws.Range("L58").Value = 279272.73
$ws.Range("M58").Value = -3825
$ws.Range("N58").Value = -279572.73
$ws.Range("H64").Value = 3558.257
$ws.Range("I64").Value = 3397.4583
$ws.Range("J64").Value = 3909.0908
$ws.Range("K64").Value = 3397.4583
$ws.Range("L64").Value = 3909.0908
$ws.Range("M64").Value = -3149.4583
$ws.Range("N64").Value = -4405.0908
$ws.Range("H67").Value = 3558.257
$ws.Range("I67").Value = 3397.4583
$ws.Range("J67").Value = 3909.0908
$ws.Range("K67").Value = 3397.4583
$ws.Range("L67").Value = 3909.0908
$ws.Range("M67").Value = -2539.4583
$ws.Range("N67").Value = -5625.0908
$ws.Range("H113").Value = 2864.4285
$ws.Range("I113").Value = 2675.6667
$ws.Range("K113").Value = 2675.6667
$ws.Range("M113").Value = 578.3332999999998
$ws.Range("H121").Value = 898.7778
$ws.Range("I121").Value = 1333.3334
$ws.Range("J121").Value = 681.5
$ws.Range("K121").Value = 4000.0002
$ws.Range("L121").Value = 2044.5
$ws.Range("M121").Value = -2253.0002
$ws.Range("N121").Value = -5538.5
$ws.Range("H137").Value = 2470.9524
$ws.Range("I137").Value = 1233.381
$ws.Range("J137").Value = 4946.095
$ws.Range("K137").Value = 3700.143
$ws.Range("L137").Value = 14838.285
$ws.Range("M137").Value = -1150.143
$ws.Range("N137").Value = -19938.285

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H35").Value = 18000
$ws.Range("I35").Value = 6000
$ws.Range("J35").Value = 20000
$ws.Range("K35").Value = 6000
$ws.Range("L35").Value = 20000
$ws.Range("M35").Value = -5594
$ws.Range("N35").Value = -20812
$ws.Range("H61").Value = 6721.92
$ws.Range("I61").Value = 5013.8335
$ws.Range("J61").Value = 11114.143
$ws.Range("K61").Value = 5013.8335
$ws.Range("L61").Value = 11114.143
$ws.Range("M61").Value = -4801.8335
$ws.Range("N61").Value = -11538.143
$ws.Range("H74").Value = 3230.1924
$ws.Range("I74").Value = 3636.0527
$ws.Range("J74").Value = 2128.5715
$ws.Range("K74").Value = 3636.0527
$ws.Range("L74").Value = 2128.5715
$ws.Range("M74").Value = -2762.0527
$ws.Range("N74").Value = -3876.5715
$ws.Range("H77").Value = 3230.1924
$ws.Range("I77").Value = 3636.0527
$ws.Range("J77").Value = 2128.5715
$ws.Range("K77").Value = 18180.2635
$ws.Range("L77").Value = 10642.8575
$ws.Range("M77").Value = -13812.2635
$ws.Range("N77").Value = -19378.8575
$ws.Range("H109").Value = 35000
$ws.Range("J109").Value = 35000
$ws.Range("L109").Value = 35000
$ws.Range("N109").Value = -37774
$ws.Range("H132").Value = 5813.2896
$ws.Range("I132").Value = 1837.6
$ws.Range("J132").Value = 8406.130999999999
$ws.Range("K132").Value = 5512.799999999999
$ws.Range("L132").Value = 25218.393
$ws.Range("M132").Value = -2982.799999999999
$ws.Range("N132").Value = -30278.393
$ws.Range("H136").Value = 6721.92
$ws.Range("I136").Value = 5013.8335
$ws.Range("J136").Value = 11114.143
$ws.Range("K136").Value = 15041.5005
$ws.Range("L136").Value = 33342.429
$ws.Range("M136").Value = -12491.5005
$ws.Range("N136").Value = -38442.429

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 9103.166999999999
$ws.Range("I7").Value = 430
$ws.Range("J7").Value = 26449.5
$ws.Range("K7").Value = 430
$ws.Range("L7").Value = 26449.5
$ws.Range("M7").Value = -317
$ws.Range("N7").Value = -26675.5
$ws.Range("H22").Value = 466.33334
$ws.Range("I22").Value = 399.5
$ws.Range("K22").Value = 399.5
$ws.Range("M22").Value = -226.5
$ws.Range("H107").Value = 2048.5217
$ws.Range("I107").Value = 2019.1904
$ws.Range("J107").Value = 2356.5
$ws.Range("K107").Value = 2019.1904
$ws.Range("L107").Value = 2356.5
$ws.Range("M107").Value = -99.19039999999995
$ws.Range("N107").Value = -6196.5
$ws.Range("H109").Value = 40684
$ws.Range("J109").Value = 40684
$ws.Range("L109").Value = 40684
$ws.Range("N109").Value = -43458

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1993.221
$ws.Range("I31").Value = 1448.3182
$ws.Range("J31").Value = 3791.4
$ws.Range("K31").Value = 1448.3182
$ws.Range("L31").Value = 3791.4
$ws.Range("M31").Value = -1153.3182
$ws.Range("N31").Value = -4381.4
$ws.Range("H34").Value = 1993.221
$ws.Range("I34").Value = 1448.3182
$ws.Range("J34").Value = 3791.4
$ws.Range("K34").Value = 1448.3182
$ws.Range("L34").Value = 3791.4
$ws.Range("M34").Value = -1246.3182
$ws.Range("N34").Value = -4195.4
$ws.Range("H58").Value = 4789566.5
$ws.Range("I58").Value = 18186514
$ws.Range("J58").Value = 4942.857
$ws.Range("K58").Value = 18186514
$ws.Range("L58").Value = 4942.857
$ws.Range("M58").Value = -18186311
$ws.Range("N58").Value = -5348.857
$ws.Range("H132").Value = 2399.923
$ws.Range("I132").Value = 2190.8635
$ws.Range("J132").Value = 3549.75
$ws.Range("K132").Value = 6572.5905
$ws.Range("L132").Value = 10649.25
$ws.Range("M132").Value = -4042.5905
$ws.Range("N132").Value = -15709.25
$ws.Range("H134").Value = 4556.486
$ws.Range("I134").Value = 4011.4
$ws.Range("K134").Value = 12034.2
$ws.Range("M134").Value = -9499.200000000001
$ws.Range("H136").Value = 4789566.5
$ws.Range("I136").Value = 18186514
$ws.Range("J136").Value = 4942.857
$ws.Range("K136").Value = 54559542
$ws.Range("L136").Value = 14828.571
$ws.Range("M136").Value = -54556992
$ws.Range("N136").Value = -19928.571

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 530.6667
$ws.Range("J33").Value = 357.4
$ws.Range("L33").Value = 2144.4
$ws.Range("N33").Value = -2710.4
$ws.Range("H56").Value = 95050.09
$ws.Range("I56").Value = 95050.09
$ws.Range("K56").Value = 95050.09
$ws.Range("M56").Value = -94520.09

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5735.8887
$ws.Range("I132").Value = 5946.4287
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 17839.2861
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -15309.2861
$ws.Range("N132").Value = -20057
$ws.Range("H134").Value = 46914.832
$ws.Range("J134").Value = 46914.832
$ws.Range("L134").Value = 140744.496
$ws.Range("N134").Value = -145814.496
$ws.Range("H136").Value = 19864.6
$ws.Range("J136").Value = 19864.6
$ws.Range("L136").Value = 59593.8
$ws.Range("N136").Value = -64693.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 789
$ws.Range("I22").Value = 600
$ws.Range("J22").Value = 812.625
$ws.Range("K22").Value = 600
$ws.Range("L22").Value = 812.625
$ws.Range("M22").Value = -305
$ws.Range("N22").Value = -1402.625
$ws.Range("H27").Value = 789
$ws.Range("I27").Value = 600
$ws.Range("J27").Value = 812.625
$ws.Range("K27").Value = 600
$ws.Range("L27").Value = 812.625
$ws.Range("M27").Value = -493
$ws.Range("N27").Value = -1026.625
$ws.Range("H132").Value = 4803.2285
$ws.Range("I132").Value = 4566.963
$ws.Range("J132").Value = 5600.625
$ws.Range("K132").Value = 13700.889
$ws.Range("L132").Value = 16801.875
$ws.Range("M132").Value = -11170.889
$ws.Range("N132").Value = -21861.875
$ws.Range("H136").Value = 6941.294
$ws.Range("I136").Value = 6093.6
$ws.Range("J136").Value = 7610.5264
$ws.Range("K136").Value = 18280.8
$ws.Range("L136").Value = 22831.5792
$ws.Range("M136").Value = -15730.8
$ws.Range("N136").Value = -27931.5792
